{"js": "// Replace the 100 math-problem texts in the single 20x5 table, cell by cell,\n// preserving each cell's run formatting (font, size) and paragraph formatting (alignment).\n// The replacements are applied in row-major reading order (row 0 col 0..4, row 1 col 0..4, ...),\n// matching the order of changes in the source diff.\nconst replacements = [[\"10+55=\",\"72-38=\"],[\"35-22=\",\"84+13=\"],[\"43+15=\",\"64+35=\"],[\"93-4=\",\"46+24=\"],[\"62+1=\",\"33+38=\"],[\"3+17=\",\"79-19=\"],[\"67-10=\",\"35-3=\"],[\"25+35=\",\"12+73=\"],[\"10+1=\",\"2+87=\"],[\"11+86=\",\"29+17=\"],[\"0+34=\",\"71+23=\"],[\"94-67=\",\"86-43=\"],[\"40+47=\",\"65-20=\"],[\"27-20=\",\"38-37=\"],[\"6+36=\",\"7+20=\"],[\"64+7=\",\"75-20=\"],[\"59-51=\",\"46+18=\"],[\"84-68=\",\"11-7=\"],[\"75-27=\",\"2+65=\"],[\"13-12=\",\"18+10=\"],[\"51+13=\",\"34-6=\"],[\"87-66=\",\"14+25=\"],[\"41+47=\",\"21+5=\"],[\"39+47=\",\"30+16=\"],[\"23+7=\",\"50+3=\"],[\"41+44=\",\"20+75=\"],[\"0+27=\",\"89-76=\"],[\"44-32=\",\"20+61=\"],[\"99-43=\",\"67+14=\"],[\"26+60=\",\"46+7=\"],[\"92-44=\",\"11+56=\"],[\"41+58=\",\"30+24=\"],[\"96-70=\",\"52-33=\"],[\"11+74=\",\"5+28=\"],[\"54-25=\",\"56+0=\"],[\"33+63=\",\"93-12=\"],[\"13+36=\",\"44+34=\"],[\"70+15=\",\"19+24=\"],[\"80-27=\",\"56-44=\"],[\"48-5=\",\"86+12=\"],[\"98-16=\",\"90-65=\"],[\"98-31=\",\"89-72=\"],[\"50-41=\",\"9+27=\"],[\"76-67=\",\"72-1=\"],[\"63+12=\",\"88-77=\"],[\"19+28=\",\"29-7=\"],[\"1+22=\",\"40-6=\"],[\"88-29=\",\"6+63=\"],[\"88-67=\",\"75-14=\"],[\"25-17=\",\"76-2=\"],[\"40+54=\",\"81-48=\"],[\"16+29=\",\"42+46=\"],[\"39+20=\",\"52+15=\"],[\"2+16=\",\"85-61=\"],[\"34+59=\",\"10+14=\"],[\"31+53=\",\"60-57=\"],[\"5+83=\",\"39+18=\"],[\"53+5=\",\"94+5=\"],[\"93+2=\",\"16+18=\"],[\"70+0=\",\"37+40=\"],[\"66-57=\",\"70-17=\"],[\"57+16=\",\"66-22=\"],[\"36+29=\",\"67+31=\"],[\"73+7=\",\"67+0=\"],[\"28+63=\",\"13+1=\"],[\"11+53=\",\"76+19=\"],[\"80-41=\",\"30-29=\"],[\"18-11=\",\"19+40=\"],[\"52-17=\",\"10+80=\"],[\"99-7=\",\"36-18=\"],[\"41+8=\",\"65+22=\"],[\"63-2=\",\"18+3=\"],[\"31+59=\",\"25+34=\"],[\"98-9=\",\"67-60=\"],[\"43+1=\",\"84-58=\"],[\"86-52=\",\"79-3=\"],[\"33-17=\",\"97-1=\"],[\"86-0=\",\"17-17=\"],[\"18+80=\",\"54+40=\"],[\"44+0=\",\"45-32=\"],[\"81-55=\",\"73-39=\"],[\"17+1=\",\"92-1=\"],[\"62-44=\",\"1+51=\"],[\"36+3=\",\"16+67=\"],[\"25+53=\",\"11+38=\"],[\"65+23=\",\"82-33=\"],[\"65+26=\",\"52-5=\"],[\"68-6=\",\"67+28=\"],[\"80-68=\",\"14+19=\"],[\"56+36=\",\"9+80=\"],[\"5-2=\",\"77-38=\"],[\"87-21=\",\"98-66=\"],[\"21+28=\",\"63+25=\"],[\"51-3=\",\"91-52=\"],[\"54-29=\",\"81+9=\"],[\"61+29=\",\"76-61=\"],[\"54+4=\",\"7+32=\"],[\"81+17=\",\"96+3=\"],[\"91-91=\",\"30+55=\"],[\"14+78=\",\"6+3=\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected at least one table in the document body.\");\n}\n\nconst table = tables.items[0];\n// NOTE: table.columnCount isn't reliably available, so derive the grid shape\n// from table.values instead (also confirms the table actually holds text).\ntable.load(\"values\");\nawait context.sync();\n\nconst rowCount = table.values.length;\nconst colCount = rowCount > 0 ? table.values[0].length : 0;\n\nif (rowCount * colCount < replacements.length) {\n  throw new Error(`Table too small: ${rowCount}x${colCount} cannot hold ${replacements.length} replacements.`);\n}\n\n// Collect the first-paragraph range of every cell (row-major order) so we can\n// batch the .load() calls before touching any text.\nconst ranges = [];\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    const idx = r * colCount + c;\n    if (idx >= replacements.length) {\n      continue;\n    }\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    const range = para.getRange();\n    range.load(\"text\");\n    ranges.push(range);\n  }\n}\nawait context.sync();\n\n// Sanity-check that the cells we collected (in row-major order) actually hold the\n// \"before\" text the diff expects at that position, so a structural mismatch fails\n// loudly instead of silently mis-writing cells.\nfor (let i = 0; i < ranges.length; i++) {\n  const expectedBefore = replacements[i][0];\n  const actual = ranges[i].text;\n  if (actual !== expectedBefore) {\n    throw new Error(\n      `Cell ${i} text mismatch: expected \"${expectedBefore}\" but found \"${actual}\".`\n    );\n  }\n}\n\n// Replace each cell's text with its \"after\" value using InsertLocation.replace so the\n// existing run's formatting (font/size) and the paragraph's formatting (alignment) are\n// kept intact.\nfor (let i = 0; i < ranges.length; i++) {\n  const after = replacements[i][1];\n  ranges[i].insertText(after, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 math-problem texts in the single 20x5 table, cell by cell,\n# preserving each cell's run formatting (font, size) and paragraph formatting\n# (alignment). Replacements are applied in row-major reading order (row 1 col\n# 1..5, row 2 col 1..5, ...), matching the order of changes in the source diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"10+55=\",\"72-38=\"),\n    @(\"35-22=\",\"84+13=\"),\n    @(\"43+15=\",\"64+35=\"),\n    @(\"93-4=\",\"46+24=\"),\n    @(\"62+1=\",\"33+38=\"),\n    @(\"3+17=\",\"79-19=\"),\n    @(\"67-10=\",\"35-3=\"),\n    @(\"25+35=\",\"12+73=\"),\n    @(\"10+1=\",\"2+87=\"),\n    @(\"11+86=\",\"29+17=\"),\n    @(\"0+34=\",\"71+23=\"),\n    @(\"94-67=\",\"86-43=\"),\n    @(\"40+47=\",\"65-20=\"),\n    @(\"27-20=\",\"38-37=\"),\n    @(\"6+36=\",\"7+20=\"),\n    @(\"64+7=\",\"75-20=\"),\n    @(\"59-51=\",\"46+18=\"),\n    @(\"84-68=\",\"11-7=\"),\n    @(\"75-27=\",\"2+65=\"),\n    @(\"13-12=\",\"18+10=\"),\n    @(\"51+13=\",\"34-6=\"),\n    @(\"87-66=\",\"14+25=\"),\n    @(\"41+47=\",\"21+5=\"),\n    @(\"39+47=\",\"30+16=\"),\n    @(\"23+7=\",\"50+3=\"),\n    @(\"41+44=\",\"20+75=\"),\n    @(\"0+27=\",\"89-76=\"),\n    @(\"44-32=\",\"20+61=\"),\n    @(\"99-43=\",\"67+14=\"),\n    @(\"26+60=\",\"46+7=\"),\n    @(\"92-44=\",\"11+56=\"),\n    @(\"41+58=\",\"30+24=\"),\n    @(\"96-70=\",\"52-33=\"),\n    @(\"11+74=\",\"5+28=\"),\n    @(\"54-25=\",\"56+0=\"),\n    @(\"33+63=\",\"93-12=\"),\n    @(\"13+36=\",\"44+34=\"),\n    @(\"70+15=\",\"19+24=\"),\n    @(\"80-27=\",\"56-44=\"),\n    @(\"48-5=\",\"86+12=\"),\n    @(\"98-16=\",\"90-65=\"),\n    @(\"98-31=\",\"89-72=\"),\n    @(\"50-41=\",\"9+27=\"),\n    @(\"76-67=\",\"72-1=\"),\n    @(\"63+12=\",\"88-77=\"),\n    @(\"19+28=\",\"29-7=\"),\n    @(\"1+22=\",\"40-6=\"),\n    @(\"88-29=\",\"6+63=\"),\n    @(\"88-67=\",\"75-14=\"),\n    @(\"25-17=\",\"76-2=\"),\n    @(\"40+54=\",\"81-48=\"),\n    @(\"16+29=\",\"42+46=\"),\n    @(\"39+20=\",\"52+15=\"),\n    @(\"2+16=\",\"85-61=\"),\n    @(\"34+59=\",\"10+14=\"),\n    @(\"31+53=\",\"60-57=\"),\n    @(\"5+83=\",\"39+18=\"),\n    @(\"53+5=\",\"94+5=\"),\n    @(\"93+2=\",\"16+18=\"),\n    @(\"70+0=\",\"37+40=\"),\n    @(\"66-57=\",\"70-17=\"),\n    @(\"57+16=\",\"66-22=\"),\n    @(\"36+29=\",\"67+31=\"),\n    @(\"73+7=\",\"67+0=\"),\n    @(\"28+63=\",\"13+1=\"),\n    @(\"11+53=\",\"76+19=\"),\n    @(\"80-41=\",\"30-29=\"),\n    @(\"18-11=\",\"19+40=\"),\n    @(\"52-17=\",\"10+80=\"),\n    @(\"99-7=\",\"36-18=\"),\n    @(\"41+8=\",\"65+22=\"),\n    @(\"63-2=\",\"18+3=\"),\n    @(\"31+59=\",\"25+34=\"),\n    @(\"98-9=\",\"67-60=\"),\n    @(\"43+1=\",\"84-58=\"),\n    @(\"86-52=\",\"79-3=\"),\n    @(\"33-17=\",\"97-1=\"),\n    @(\"86-0=\",\"17-17=\"),\n    @(\"18+80=\",\"54+40=\"),\n    @(\"44+0=\",\"45-32=\"),\n    @(\"81-55=\",\"73-39=\"),\n    @(\"17+1=\",\"92-1=\"),\n    @(\"62-44=\",\"1+51=\"),\n    @(\"36+3=\",\"16+67=\"),\n    @(\"25+53=\",\"11+38=\"),\n    @(\"65+23=\",\"82-33=\"),\n    @(\"65+26=\",\"52-5=\"),\n    @(\"68-6=\",\"67+28=\"),\n    @(\"80-68=\",\"14+19=\"),\n    @(\"56+36=\",\"9+80=\"),\n    @(\"5-2=\",\"77-38=\"),\n    @(\"87-21=\",\"98-66=\"),\n    @(\"21+28=\",\"63+25=\"),\n    @(\"51-3=\",\"91-52=\"),\n    @(\"54-29=\",\"81+9=\"),\n    @(\"61+29=\",\"76-61=\"),\n    @(\"54+4=\",\"7+32=\"),\n    @(\"81+17=\",\"96+3=\"),\n    @(\"91-91=\",\"30+55=\"),\n    @(\"14+78=\",\"6+3=\")\n)\n\nif ($d.Tables.Count -lt 1) {\n    throw \"Expected at least one table in the document.\"\n}\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nif (($rowCount * $colCount) -lt $replacements.Count) {\n    throw \"Table too small: $rowCount x $colCount cannot hold $($replacements.Count) replacements.\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $replacements.Count) {\n            break\n        }\n\n        $before = $replacements[$idx][0]\n        $after = $replacements[$idx][1]\n\n        $cell = $table.Cell($r, $c)\n        $find = $cell.Range.Find\n        $find.ClearFormatting()\n        $find.Text = $before\n        $find.Replacement.ClearFormatting()\n        $find.Replacement.Text = $after\n        # Scope the search to this single cell's Range and replace only the\n        # first (and only) match so cell text that happens to be a substring\n        # of another cell's text (e.g. \"6+36=\" inside \"56+36=\") is never\n        # touched by mistake.\n        $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n        if (-not $found) {\n            throw \"Cell ($r,$c) text mismatch: expected to find `\"$before`\" but it was not present.\"\n        }\n\n        $idx++\n    }\n}\n\nWrite-Output \"Replaced $idx cell(s).\"\n"}
